$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 470.2439
$ws.Range("J17").Value = 465.425
$ws.Range("L17").Value = 1396.275
$ws.Range("N17").Value = -1732.275
$ws.Range("H20").Value = 170
$ws.Range("I20").Value = 170
$ws.Range("K20").Value = 170
$ws.Range("M20").Value = 60
$ws.Range("H35").Value = 170
$ws.Range("I35").Value = 170
$ws.Range("K35").Value = 170
$ws.Range("M35").Value = 209
$ws.Range("H40").Value = 1087.35
$ws.Range("J40").Value = 1371
$ws.Range("L40").Value = 1371
$ws.Range("N40").Value = -1721
$ws.Range("H53").Value = 2555.923
$ws.Range("I53").Value = 262
$ws.Range("J53").Value = 3575.4443
$ws.Range("K53").Value = 262
$ws.Range("L53").Value = 3575.4443
$ws.Range("M53").Value = 375
$ws.Range("N53").Value = -4849.4443
$ws.Range("H74").Value = 2355.1035
$ws.Range("I74").Value = 1645.15
$ws.Range("J74").Value = 3932.7778
$ws.Range("K74").Value = 1645.15
$ws.Range("L74").Value = 3932.7778
$ws.Range("M74").Value = -709.1500000000001
$ws.Range("N74").Value = -5804.7778
$ws.Range("H77").Value = 2355.1035
$ws.Range("I77").Value = 1645.15
$ws.Range("J77").Value = 3932.7778
$ws.Range("K77").Value = 8225.75
$ws.Range("L77").Value = 19663.889
$ws.Range("M77").Value = -3545.75
$ws.Range("N77").Value = -29023.889
$ws.Range("H80").Value = 8086.1113
$ws.Range("J80").Value = 10523.077
$ws.Range("L80").Value = 31569.231
$ws.Range("N80").Value = -33565.231
$ws.Range("H83").Value = 8086.1113
$ws.Range("J83").Value = 10523.077
$ws.Range("L83").Value = 94707.693
$ws.Range("N83").Value = -104691.693
$ws.Range("H88").Value = 1654.1111
$ws.Range("I88").Value = 750
$ws.Range("K88").Value = 750
$ws.Range("M88").Value = -344
$ws.Range("H91").Value = 1654.1111
$ws.Range("I91").Value = 750
$ws.Range("K91").Value = 750
$ws.Range("M91").Value = 654
$ws.Range("H138").Value = 2502.383
$ws.Range("J138").Value = 3331.9395
$ws.Range("L138").Value = 9995.818499999999
$ws.Range("N138").Value = -20275.8185

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 4248.933
$ws.Range("I102").Value = 3778
$ws.Range("J102").Value = 4366.6665
$ws.Range("K102").Value = 3778
$ws.Range("L102").Value = 4366.6665
$ws.Range("M102").Value = -2156
$ws.Range("N102").Value = -7610.6665

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 18396.727
$ws.Range("J82").Value = 49811.668
$ws.Range("L82").Value = 49811.668
$ws.Range("N82").Value = -50577.668
$ws.Range("H85").Value = 18396.727
$ws.Range("J85").Value = 49811.668
$ws.Range("L85").Value = 49811.668
$ws.Range("N85").Value = -52463.668
$ws.Range("H86").Value = 40200.54
$ws.Range("I86").Value = 46736.363
$ws.Range("J86").Value = 4253.5
$ws.Range("K86").Value = 46736.363
$ws.Range("L86").Value = 4253.5
$ws.Range("M86").Value = -45613.363
$ws.Range("N86").Value = -6499.5
$ws.Range("H89").Value = 40200.54
$ws.Range("I89").Value = 46736.363
$ws.Range("J89").Value = 4253.5
$ws.Range("K89").Value = 233681.815
$ws.Range("L89").Value = 21267.5
$ws.Range("M89").Value = -228065.815
$ws.Range("N89").Value = -32499.5
$ws.Range("H94").Value = 4084.5881
$ws.Range("I94").Value = 1306
$ws.Range("J94").Value = 6554.4443
$ws.Range("K94").Value = 1306
$ws.Range("L94").Value = 6554.4443
$ws.Range("M94").Value = -855
$ws.Range("N94").Value = -7456.4443
$ws.Range("H99").Value = 955
$ws.Range("I99").Value = 947
$ws.Range("J99").Value = 995
$ws.Range("K99").Value = 947
$ws.Range("L99").Value = 995
$ws.Range("M99").Value = 551
$ws.Range("N99").Value = -3991
$ws.Range("H105").Value = 1899.9286
$ws.Range("I105").Value = 1714.1428
$ws.Range("K105").Value = 1714.1428
$ws.Range("M105").Value = 32.85719999999992
$ws.Range("H107").Value = 702.03705
$ws.Range("I107").Value = 740.6818
$ws.Range("J107").Value = 532
$ws.Range("K107").Value = 740.6818
$ws.Range("L107").Value = 532
$ws.Range("M107").Value = 1179.3182
$ws.Range("N107").Value = -4372

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11614.5
$ws.Range("I31").Value = 26053.924
$ws.Range("J31").Value = 2675.8096
$ws.Range("K31").Value = 26053.924
$ws.Range("L31").Value = 2675.8096
$ws.Range("M31").Value = -25758.924
$ws.Range("N31").Value = -3265.8096
$ws.Range("H34").Value = 11614.5
$ws.Range("I34").Value = 26053.924
$ws.Range("J34").Value = 2675.8096
$ws.Range("K34").Value = 26053.924
$ws.Range("L34").Value = 2675.8096
$ws.Range("M34").Value = -25851.924
$ws.Range("N34").Value = -3079.8096

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 628.075
$ws.Range("J5").Value = 823.8421
$ws.Range("L5").Value = 2471.5263
$ws.Range("N5").Value = -2695.5263
$ws.Range("H131").Value = 124289.96
$ws.Range("J131").Value = 134162.58
$ws.Range("L131").Value = 402487.74
$ws.Range("N131").Value = -412567.74
$ws.Range("H133").Value = 3153
$ws.Range("I133").Value = 2547.5
$ws.Range("J133").Value = 3960.3333
$ws.Range("K133").Value = 7642.5
$ws.Range("L133").Value = 11880.9999
$ws.Range("M133").Value = -2582.5
$ws.Range("N133").Value = -22000.9999
$ws.Range("H134").Value = 6634.4287
$ws.Range("I134").Value = 7043.278
$ws.Range("K134").Value = 21129.834
$ws.Range("M134").Value = -16059.834
$ws.Range("H135").Value = 628.075
$ws.Range("J135").Value = 823.8421
$ws.Range("L135").Value = 7414.5789
$ws.Range("N135").Value = -12484.5789

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 7739.476
$ws.Range("I80").Value = 10329.083
$ws.Range("J80").Value = 4286.6665
$ws.Range("K80").Value = 10329.083
$ws.Range("L80").Value = 4286.6665
$ws.Range("M80").Value = -9331.083000000001
$ws.Range("N80").Value = -6282.6665
$ws.Range("H83").Value = 7739.476
$ws.Range("I83").Value = 10329.083
$ws.Range("J83").Value = 4286.6665
$ws.Range("K83").Value = 51645.415
$ws.Range("L83").Value = 21433.3325
$ws.Range("M83").Value = -46653.415
$ws.Range("N83").Value = -31417.3325
$ws.Range("H132").Value = 128199.414
$ws.Range("I132").Value = 202560.2
$ws.Range("K132").Value = 607680.6000000001
$ws.Range("M132").Value = -605150.6000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4299.4
$ws.Range("I40").Value = 3856.2856
$ws.Range("J40").Value = 5333.3335
$ws.Range("K40").Value = 3856.2856
$ws.Range("L40").Value = 5333.3335
$ws.Range("M40").Value = -3720.2856
$ws.Range("N40").Value = -5605.3335
$ws.Range("H61").Value = 4567.3335
$ws.Range("I61").Value = 1625.3334
$ws.Range("J61").Value = 8490
$ws.Range("K61").Value = 1625.3334
$ws.Range("L61").Value = 8490
$ws.Range("M61").Value = -1423.3334
$ws.Range("N61").Value = -8894
$ws.Range("H100").Value = 5692.1665
$ws.Range("I100").Value = 1038.25
$ws.Range("K100").Value = 1038.25
$ws.Range("M100").Value = -497.25
$ws.Range("H113").Value = 4567.3335
$ws.Range("I113").Value = 1625.3334
$ws.Range("J113").Value = 8490
$ws.Range("K113").Value = 1625.3334
$ws.Range("L113").Value = 8490
$ws.Range("M113").Value = 544.6666
$ws.Range("N113").Value = -12830
$ws.Range("H122").Value = 3646.2727
$ws.Range("I122").Value = 3526
$ws.Range("J122").Value = 3715
$ws.Range("K122").Value = 10578
$ws.Range("L122").Value = 11145
$ws.Range("M122").Value = -8128
$ws.Range("N122").Value = -16045

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1373.5714
$ws.Range("I81").Value = 1522.6
$ws.Range("J81").Value = 1001
$ws.Range("K81").Value = 3045.2
$ws.Range("L81").Value = 2002
$ws.Range("M81").Value = -1984.2
$ws.Range("N81").Value = -4124
$ws.Range("H84").Value = 1373.5714
$ws.Range("I84").Value = 1522.6
$ws.Range("J84").Value = 1001
$ws.Range("K84").Value = 15226
$ws.Range("L84").Value = 10010
$ws.Range("M84").Value = -9922
$ws.Range("N84").Value = -20618
$ws.Range("H96").Value = 3976.8667
$ws.Range("J96").Value = 4895
$ws.Range("L96").Value = 4895
$ws.Range("N96").Value = -7641
